$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.647790253162384
$ws.Range("B1").Value = 0.7441759705543518
$ws.Range("C1").Value = 4.895464897155762
$ws.Range("D1").Value = 1.896253705024719
$ws.Range("E1").Value = 0.7260017395019531
